$d = $word.ActiveDocument

# The paragraph currently reads "Versi" + "on" + " 2" + "." (split across
# several runs). We want the visible text "Version 1." while keeping the
# proofErr spell-check markers around "Version" and the _GoBack bookmark
# in place, and dropping the now-separate "." run.

# Step 1: merge the "Versi"/"on" runs into a single "Version" run.
# (Delete "on" so it collapses into the "Versi" run, then re-insert "on"
# right after it so both words live in the same run.)
$r1 = $d.Range(5, 7)
$r1.Text = ""
$r2 = $d.Range(5, 5)
$r2.InsertAfter("on")

# Step 2: change the "2" run's text to "1." (stays in its own run).
$r3 = $d.Range(8, 9)
$r3.Text = "1."

# Step 3: remove the trailing "." run that followed the bookmark.
$r4 = $d.Range(10, 11)
$r4.Text = ""
